# Quarterly indexing esoteric bug-fix operation
#
# Column A holds quarterly "as-of" dates that were incorrectly anchored to
# the 1st of a quarter month. The fix re-anchors every date in column A to
# the 15th of the month following the original date, i.e.
#   new_date = 15th of the month that is 1 month after old_date
# Since every existing value in column A is the 1st of a month,
# EDATE(old_date, 1) lands on the 1st of the target month, so adding 14
# days gives the 15th of that month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$wf = $excel.WorksheetFunction

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldSerial = $cell.Value2
    if ($oldSerial -ne $null) {
        $nextMonth1st = $wf.EDate($oldSerial, 1)
        $newSerial = $nextMonth1st + 14
        $cell.Value = $newSerial
    }
}
